$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for b.md (row 3) reflects the new handoff status/date ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-07 07:59:56"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
# Error Detail column widened to fit the new message text
$zhcn.Columns.Item(16).ColumnWidth = 39.17
# a.md status -> Ready for handoff
$zhcn.Range("C2").Value = "Ready for handoff"
# b.md row: status, duplicate flag, latest handoff file/datetime, error detail
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-07 07:59:44"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/53ab17dca758528764e1a6b465bdf928c8ca28ff/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8192a6605f87a9913d5f52057fe41e5b49342f1b/e2e/b.md."

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
# Error Detail column widened to fit the new message text
$dede.Columns.Item(16).ColumnWidth = 39.17
# a.md status -> Ready for handoff
$dede.Range("C2").Value = "Ready for handoff"
# b.md row: status, duplicate flag, latest handoff file/datetime, error detail
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-07 07:59:56"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/53ab17dca758528764e1a6b465bdf928c8ca28ff/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8192a6605f87a9913d5f52057fe41e5b49342f1b/e2e/b.md."
